$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 16 new data rows (11-26) to the "nueva bobina" log.
#
# The sheet stores every column as text (matching the existing rows 2-10),
# so each row is written with NumberFormat "@" applied first to avoid
# Excel's automatic text->number coercion, then the style is reset back to
# "Normal" (the sheet's default/only style) so no extra formatting is left
# on the cells. The "Sec" column (F) mirrors the source data exactly: most
# entries are text, but F11 and F13 are real numbers (matching the existing
# F2/F4/F9/F10 pattern already present in the sheet), so those two cells get
# their value reassigned as a number after the style reset.

$rng = $ws.Range("A11:K11")
$rng.NumberFormat = "@"
$ws.Range("A11").Value = '205'
$ws.Range("B11").Value = '120'
$ws.Range("C11").Value = '150'
$ws.Range("D11").Value = '1771'
$ws.Range("E11").Value = '8205'
$ws.Range("F11").Value = '1'
$ws.Range("G11").Value = '86337'
$ws.Range("H11").Value = '2025-07-01 14:29'
$ws.Range("I11").Value = 'A'
$ws.Range("J11").Value = '03'
$ws.Range("K11").Value = 'L.BLANCO'
$rng.Style = "Normal"
$ws.Range("F11").Value = 1

$rng = $ws.Range("A12:K12")
$rng.NumberFormat = "@"
$ws.Range("A12").Value = '35'
$ws.Range("B12").Value = '120'
$ws.Range("C12").Value = '150'
$ws.Range("D12").Value = '296'
$ws.Range("E12").Value = '8205'
$ws.Range("F12").Value = '2'
$ws.Range("G12").Value = '86337'
$ws.Range("H12").Value = '2025-07-01 14:29'
$ws.Range("I12").Value = 'A'
$ws.Range("J12").Value = '03'
$ws.Range("K12").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A13:K13")
$rng.NumberFormat = "@"
$ws.Range("A13").Value = '205'
$ws.Range("B13").Value = '120'
$ws.Range("C13").Value = '150'
$ws.Range("D13").Value = '1763'
$ws.Range("E13").Value = '8207'
$ws.Range("F13").Value = '1'
$ws.Range("G13").Value = '86337'
$ws.Range("H13").Value = '2025-07-01 14:31'
$ws.Range("I13").Value = 'A'
$ws.Range("J13").Value = '03'
$ws.Range("K13").Value = 'L.BLANCO'
$rng.Style = "Normal"
$ws.Range("F13").Value = 1

$rng = $ws.Range("A14:K14")
$rng.NumberFormat = "@"
$ws.Range("A14").Value = '35'
$ws.Range("B14").Value = '120'
$ws.Range("C14").Value = '150'
$ws.Range("D14").Value = '292'
$ws.Range("E14").Value = '8207'
$ws.Range("F14").Value = '2'
$ws.Range("G14").Value = '86337'
$ws.Range("H14").Value = '2025-07-01 14:31'
$ws.Range("I14").Value = 'A'
$ws.Range("J14").Value = '03'
$ws.Range("K14").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A15:K15")
$rng.NumberFormat = "@"
$ws.Range("A15").Value = '205'
$ws.Range("B15").Value = '120'
$ws.Range("C15").Value = '150'
$ws.Range("D15").Value = '1760'
$ws.Range("E15").Value = '8209'
$ws.Range("F15").Value = '1'
$ws.Range("G15").Value = '86337'
$ws.Range("H15").Value = '2025-07-01 14:31'
$ws.Range("I15").Value = 'A'
$ws.Range("J15").Value = '03'
$ws.Range("K15").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A16:K16")
$rng.NumberFormat = "@"
$ws.Range("A16").Value = '35'
$ws.Range("B16").Value = '120'
$ws.Range("C16").Value = '150'
$ws.Range("D16").Value = '295'
$ws.Range("E16").Value = '8209'
$ws.Range("F16").Value = '2'
$ws.Range("G16").Value = '86337'
$ws.Range("H16").Value = '2025-07-01 14:31'
$ws.Range("I16").Value = 'A'
$ws.Range("J16").Value = '03'
$ws.Range("K16").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A17:K17")
$rng.NumberFormat = "@"
$ws.Range("A17").Value = '205'
$ws.Range("B17").Value = '120'
$ws.Range("C17").Value = '150'
$ws.Range("D17").Value = '1801'
$ws.Range("E17").Value = '8211'
$ws.Range("F17").Value = '1'
$ws.Range("G17").Value = '86337'
$ws.Range("H17").Value = '2025-07-01 14:31'
$ws.Range("I17").Value = 'A'
$ws.Range("J17").Value = '03'
$ws.Range("K17").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A18:K18")
$rng.NumberFormat = "@"
$ws.Range("A18").Value = '35'
$ws.Range("B18").Value = '120'
$ws.Range("C18").Value = '150'
$ws.Range("D18").Value = '299'
$ws.Range("E18").Value = '8211'
$ws.Range("F18").Value = '2'
$ws.Range("G18").Value = '86337'
$ws.Range("H18").Value = '2025-07-01 14:31'
$ws.Range("I18").Value = 'A'
$ws.Range("J18").Value = '03'
$ws.Range("K18").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A19:K19")
$rng.NumberFormat = "@"
$ws.Range("A19").Value = '120'
$ws.Range("B19").Value = '120'
$ws.Range("C19").Value = '150'
$ws.Range("D19").Value = '555'
$ws.Range("E19").Value = '3323'
$ws.Range("F19").Value = '1'
$ws.Range("G19").Value = '76767'
$ws.Range("H19").Value = '2025-07-10 10:16'
$ws.Range("I19").Value = 'B'
$ws.Range("J19").Value = '03'
$ws.Range("K19").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A20:K20")
$rng.NumberFormat = "@"
$ws.Range("A20").Value = '110'
$ws.Range("B20").Value = '120'
$ws.Range("C20").Value = '150'
$ws.Range("D20").Value = '436'
$ws.Range("E20").Value = '3323'
$ws.Range("F20").Value = '2'
$ws.Range("G20").Value = '76767'
$ws.Range("H20").Value = '2025-07-10 10:16'
$ws.Range("I20").Value = 'B'
$ws.Range("J20").Value = '03'
$ws.Range("K20").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A21:K21")
$rng.NumberFormat = "@"
$ws.Range("A21").Value = '120'
$ws.Range("B21").Value = '120'
$ws.Range("C21").Value = '150'
$ws.Range("D21").Value = '578'
$ws.Range("E21").Value = '3324'
$ws.Range("F21").Value = '1'
$ws.Range("G21").Value = '76767'
$ws.Range("H21").Value = '2025-07-10 10:16'
$ws.Range("I21").Value = 'B'
$ws.Range("J21").Value = '03'
$ws.Range("K21").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A22:K22")
$rng.NumberFormat = "@"
$ws.Range("A22").Value = '110'
$ws.Range("B22").Value = '120'
$ws.Range("C22").Value = '150'
$ws.Range("D22").Value = '454'
$ws.Range("E22").Value = '3324'
$ws.Range("F22").Value = '2'
$ws.Range("G22").Value = '76767'
$ws.Range("H22").Value = '2025-07-10 10:16'
$ws.Range("I22").Value = 'B'
$ws.Range("J22").Value = '03'
$ws.Range("K22").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A23:K23")
$rng.NumberFormat = "@"
$ws.Range("A23").Value = '120'
$ws.Range("B23").Value = '120'
$ws.Range("C23").Value = '150'
$ws.Range("D23").Value = '577'
$ws.Range("E23").Value = '3325'
$ws.Range("F23").Value = '1'
$ws.Range("G23").Value = '76767'
$ws.Range("H23").Value = '2025-07-10 10:16'
$ws.Range("I23").Value = 'B'
$ws.Range("J23").Value = '03'
$ws.Range("K23").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A24:K24")
$rng.NumberFormat = "@"
$ws.Range("A24").Value = '110'
$ws.Range("B24").Value = '120'
$ws.Range("C24").Value = '150'
$ws.Range("D24").Value = '434'
$ws.Range("E24").Value = '3325'
$ws.Range("F24").Value = '2'
$ws.Range("G24").Value = '76767'
$ws.Range("H24").Value = '2025-07-10 10:16'
$ws.Range("I24").Value = 'B'
$ws.Range("J24").Value = '03'
$ws.Range("K24").Value = 'L.BLANCO'
$rng.Style = "Normal"

$rng = $ws.Range("A25:K25")
$rng.NumberFormat = "@"
$ws.Range("A25").Value = '80'
$ws.Range("B25").Value = '120'
$ws.Range("C25").Value = '130'
$ws.Range("D25").Value = '540'
$ws.Range("E25").Value = '1278'
$ws.Range("F25").Value = '1'
$ws.Range("G25").Value = '879'
$ws.Range("H25").Value = '2025-07-28 16:56'
$ws.Range("I25").Value = 'A'
$ws.Range("J25").Value = '06'
$ws.Range("K25").Value = 'LINER PER'
$rng.Style = "Normal"

$rng = $ws.Range("A26:K26")
$rng.NumberFormat = "@"
$ws.Range("A26").Value = '90'
$ws.Range("B26").Value = '120'
$ws.Range("C26").Value = '130'
$ws.Range("D26").Value = '654'
$ws.Range("E26").Value = '1278'
$ws.Range("F26").Value = '2'
$ws.Range("G26").Value = '879'
$ws.Range("H26").Value = '2025-07-28 16:56'
$ws.Range("I26").Value = 'A'
$ws.Range("J26").Value = '06'
$ws.Range("K26").Value = 'LINER PER'
$rng.Style = "Normal"
